# Update "Data availability index" (column G) values for several rows
# as part of addressing the spatial coverage issue.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G42").Value = 75.41
$ws.Range("G44").Value = 60.72
$ws.Range("G46").Value = 66.81999999999999
$ws.Range("G47").Value = 63.73
$ws.Range("G50").Value = 62.31
$ws.Range("G61").Value = 63.09
$ws.Range("G75").Value = 60.53
$ws.Range("G83").Value = 62.89
